$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.481.54"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "1.836.75"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.51"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5371"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2984"
$ws.Range("E8").Value = "  -7.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06927"
$ws.Range("E9").Value = "  +1.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.52"
$ws.Range("E10").Value = "  -7.74%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.844.69"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7315"
$ws.Range("E12").Value = "  -6.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07204"
$ws.Range("E13").Value = "  -7.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.12"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.985"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.81"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007894"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "26.499.80"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").Value = "2.076.80"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.573"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.985"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.186"
$ws.Range("E24").Value = "  -3.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.65"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.166"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.717"
$ws.Range("E27").Value = "  +2.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.96"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.90"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.225"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08844"
$ws.Range("E31").Value = "  +1.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.037"
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04837"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.934"
$ws.Range("E34").Value = "  +2.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7225"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.130"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.092"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.282"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01713"
$ws.Range("E39").Value = "  -4.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4701"
$ws.Range("E40").Value = "  -3.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9033"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "107.67"
$ws.Range("E42").Value = "  -3.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.889"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.390"
$ws.Range("E45").Value = "  -3.85%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.046"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1246"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4068"
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.77"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05765"
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8914"
$ws.Range("E51").Value = "  +0.09%  "
